$d = $word.ActiveDocument

# The contact-info paragraph is the 2nd paragraph in the document.
$para = $d.Paragraphs(2)
$rng = $para.Range

# Narrow the range to just the leading text before the e-mail hyperlink
# (i.e. "3/27 Mahia Tce  Kings Beach, QLD 4551  DOB: 16/10/1981  0458385114  ").
$full = $rng.Text
$cut = $full.IndexOf("comital333")
$rng.End = $rng.Start + $cut

$rng.Text = " 0458385114  "
